$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '66.589.98'
Set-TextValue 'E2' '  -1.84%  '
Set-TextValue 'D3' '3.488.05'
Set-TextValue 'E3' '  -0.34%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '583.25'
Set-TextValue 'E5' '  -2.63%  '
Set-TextValue 'D6' '174.29'
Set-TextValue 'E6' '  -3.80%  '
Set-TextValue 'D7' '0.999'
Set-TextValue 'E7' '  +0.04%  '
Set-TextValue 'D8' '3.490.81'
Set-TextValue 'E8' '  -0.28%  '
Set-TextValue 'E9' '  -3.57%  '
Set-TextValue 'D10' '0.133'
Set-TextValue 'E10' '  -5.31%  '
Set-TextValue 'D11' '6.90'
Set-TextValue 'E11' '  -1.95%  '
Set-TextValue 'D12' '0.421'
Set-TextValue 'E12' '  -3.69%  '
Set-TextValue 'D13' '4.094.59'
Set-TextValue 'E13' '  -0.12%  '
Set-TextValue 'D14' '30.27'
Set-TextValue 'E14' '  -5.66%  '
Set-TextValue 'D16' '66.362.33'
Set-TextValue 'E16' '  -2.13%  '
Set-TextValue 'D17' '0.0000173'
Set-TextValue 'E17' '  -3.59%  '
Set-TextValue 'D18' '3.488.50'
Set-TextValue 'E18' '  -0.10%  '
Set-TextValue 'D19' '6.01'
Set-TextValue 'E19' '  -5.27%  '
Set-TextValue 'D20' '13.84'
Set-TextValue 'E20' '  -3.37%  '
Set-TextValue 'D21' '380.40'
Set-TextValue 'E21' '  -3.15%  '
Set-TextValue 'D22' '7.83'
Set-TextValue 'E22' '  -2.00%  '
Set-TextValue 'D23' '0.548'
Set-TextValue 'E23' '  +1.01%  '
Set-TextValue 'E24' '  +0.08%  '
Set-TextValue 'B25' 'Litecoin'
Set-TextValue 'C25' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D25' '72.12'
Set-TextValue 'E25' '  -1.39%  '
Set-TextValue 'B26' 'LEO'
Set-TextValue 'C26' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D26' '5.73'
Set-TextValue 'E26' '  -0.43%  '
Set-TextValue 'D27' '0.0000120'
Set-TextValue 'E27' '  -3.24%  '
Set-TextValue 'D28' '9.84'
Set-TextValue 'E28' '  -5.24%  '
Set-TextValue 'E29' '  -1.39%  '
Set-TextValue 'E30' '  +0.05%  '
Set-TextValue 'D31' '24.40'
Set-TextValue 'E31' '  +3.13%  '
Set-TextValue 'D32' '5.87'
Set-TextValue 'E32' '  -5.23%  '
Set-TextValue 'E33' '  -3.46%  '
Set-TextValue 'B34' 'Fetch.AI'
Set-TextValue 'C34' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D34' '1.32'
Set-TextValue 'E34' '  -7.46%  '
Set-TextValue 'B35' 'USDe'
Set-TextValue 'C35' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D35' '0.999'
Set-TextValue 'E35' '  -0.10%  '
Set-TextValue 'D36' '7.15'
Set-TextValue 'E36' '  -3.98%  '
Set-TextValue 'D37' '1.58'
Set-TextValue 'E37' '  -2.51%  '
Set-TextValue 'B38' 'EnergySwap'
Set-TextValue 'C38' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D38' '29.48'
Set-TextValue 'E38' '  +11.77%  '
Set-TextValue 'B39' 'Monero'
Set-TextValue 'C39' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D39' '159.54'
Set-TextValue 'E39' '  -1.85%  '
Set-TextValue 'B40' 'Mantle'
Set-TextValue 'C40' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D40' '0.893'
Set-TextValue 'E40' '  +0.54%  '
Set-TextValue 'D41' '1.78'
Set-TextValue 'E41' '  -6.07%  '
Set-TextValue 'E42' '  -3.02%  '
Set-TextValue 'B43' 'dogwifhat'
Set-TextValue 'C43' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D43' '2.54'
Set-TextValue 'E43' '  -11.00%  '
Set-TextValue 'B44' 'RenderToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D44' '6.41'
Set-TextValue 'E44' '  -6.17%  '
Set-TextValue 'B45' 'Hedera'
Set-TextValue 'C45' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D45' '0.0696'
Set-TextValue 'E45' '  -4.46%  '
Set-TextValue 'B46' 'Maker'
Set-TextValue 'C46' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D46' '2.676.12'
Set-TextValue 'E46' '  -5.46%  '
Set-TextValue 'D47' '40.63'
Set-TextValue 'E47' '  -2.28%  '
Set-TextValue 'D48' '24.51'
Set-TextValue 'E48' '  -9.00%  '
Set-TextValue 'D49' '0.0291'
Set-TextValue 'E49' '  -3.77%  '
Set-TextValue 'D50' '315.67'
Set-TextValue 'E50' '  -5.89%  '
Set-TextValue 'E51' '  -4.70%  '
